# Update data values in row 5 (Sheet1)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = 2
$ws.Range("F5").Value = -3
$ws.Range("H5").Value = 46

# Update the selected cell to reflect the author's last cursor position
$ws.Range("C5").Select()
